$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores plain-looking decimals ("360.20", "3.10", ...) as TEXT
# in the source sheet (t="inlineStr"). Assigning those strings straight to
# .Value lets Excel reinterpret them as numbers and silently drop
# significant trailing zeros (e.g. "360.20" -> 360.2), so a leading
# apostrophe is included to force text entry, exactly like typing
# '360.20 into a cell by hand in Excel.

$ws.Range("D2").Value = "'51.989.35"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "'2.791.54"
$ws.Range("E3").Value = "  -1.80%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'360.20"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "'109.65"
$ws.Range("E6").Value = "  -3.49%  "

$ws.Range("E7").Value = "  -2.76%  "

$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  -2.76%  "

$ws.Range("D10").Value = "'40.23"
$ws.Range("E10").Value = "  -3.37%  "

$ws.Range("D11").Value = "'0.0851"
$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("E13").Value = "  -2.82%  "

$ws.Range("D14").Value = "'7.58"
$ws.Range("E14").Value = "  -3.20%  "

$ws.Range("D15").Value = "'3.234.27"

$ws.Range("D16").Value = "'2.796.22"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").Value = "'0.942"
$ws.Range("E17").Value = "  +3.79%  "

$ws.Range("D18").Value = "'51.935.74"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").Value = "'3.10"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").Value = "'13.17"
$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("D22").Value = "'0.0₃0976"
$ws.Range("E22").Value = "  -1.97%  "

$ws.Range("D23").Value = "'70.38"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'270.97"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  -2.23%  "

$ws.Range("D26").Value = "'26.55"
$ws.Range("E26").Value = "  -2.30%  "

$ws.Range("D28").Value = "'0.161"
$ws.Range("E28").Value = "  +14.63%  "

$ws.Range("D29").Value = "'10.31"
$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("D30").Value = "'2.30"
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").Value = "'51.87"
$ws.Range("E32").Value = "  -3.72%  "

$ws.Range("D33").Value = "'34.34"
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("D34").Value = "'5.76"
$ws.Range("E34").Value = "  -2.40%  "

$ws.Range("D35").Value = "'0.0848"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").Value = "'5.23"
$ws.Range("E36").Value = "  -3.08%  "

$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("E38").Value = "  +3.28%  "

$ws.Range("D39").Value = "'3.21"
$ws.Range("E39").Value = "  -2.03%  "

$ws.Range("E40").Value = "  -3.95%  "

$ws.Range("D41").Value = "'2.64"
$ws.Range("E41").Value = "  +2.99%  "

$ws.Range("E42").Value = "  -2.20%  "

$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("D44").Value = "'119.55"
$ws.Range("E44").Value = "  -6.93%  "

$ws.Range("D45").Value = "'21.98"
$ws.Range("E45").Value = "  -7.81%  "

$ws.Range("D46").Value = "'2.082.94"
$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").Value = "'3.27"
$ws.Range("E47").Value = "  -4.72%  "

$ws.Range("D49").Value = "'5.83"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("E50").Value = "  -5.85%  "

$ws.Range("E51").Value = "  -2.99%  "
